$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a numeric-looking string (e.g. "1.00", "23.50")
# that must remain literal text (matching the source inline-string cell).
# Pre-format as Text so Excel does not silently coerce it to a number and
# drop meaningful trailing zeros / introduce float rounding noise.
$textCells = @(
    "D5",
    "D6",
    "D9",
    "D10",
    "D17",
    "D19",
    "D20",
    "D21",
    "D22",
    "D23",
    "D26",
    "D27",
    "D28",
    "D29",
    "D30",
    "D31",
    "D36",
    "D37",
    "D38",
    "D39",
    "D40",
    "D42",
    "D43",
    "D44",
    "D45",
    "D47",
    "D48",
    "D49",
    "D50",
    "D51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cryptos snapshot values
$ws.Range("D2").Value = "61.477.56"
$ws.Range("D3").Value = "3.378.41"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "576.76"
$ws.Range("E5").Value = "  +0.10%  "
$ws.Range("D6").Value = "136.00"
$ws.Range("E6").Value = "  +7.61%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "3.377.42"
$ws.Range("E8").Value = "  -0.79%  "
$ws.Range("D9").Value = "0.478"
$ws.Range("E9").Value = "  +0.36%  "
$ws.Range("D10").Value = "7.56"
$ws.Range("E10").Value = "  +2.27%  "
$ws.Range("E11").Value = "  +1.84%  "
$ws.Range("E12").Value = "  +1.84%  "
$ws.Range("D13").Value = "3.955.37"
$ws.Range("E13").Value = "  -0.77%  "
$ws.Range("E14").Value = "  +1.08%  "
$ws.Range("E15").Value = "  +0.98%  "
$ws.Range("D16").Value = "3.378.20"
$ws.Range("E16").Value = "  -0.85%  "
$ws.Range("D17").Value = "25.26"
$ws.Range("E17").Value = "  +1.47%  "
$ws.Range("D18").Value = "61.517.68"
$ws.Range("E18").Value = "  -2.24%  "
$ws.Range("D19").Value = "14.03"
$ws.Range("E19").Value = "  +6.18%  "
$ws.Range("D20").Value = "5.81"
$ws.Range("E20").Value = "  +2.13%  "
$ws.Range("D21").Value = "9.36"
$ws.Range("E21").Value = "  -2.29%  "
$ws.Range("D22").Value = "385.17"
$ws.Range("E22").Value = "  +2.01%  "
$ws.Range("D23").Value = "0.568"
$ws.Range("E23").Value = "  +1.51%  "
$ws.Range("D24").Value = "3.512.11"
$ws.Range("E24").Value = "  -0.89%  "
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("D26").Value = "71.01"
$ws.Range("E26").Value = "  -2.18%  "
$ws.Range("D27").Value = "0.0000119"
$ws.Range("E27").Value = "  +9.71%  "
$ws.Range("D28").Value = "1.70"
$ws.Range("E28").Value = "  +21.16%  "
$ws.Range("D29").Value = "7.84"
$ws.Range("E29").Value = "  +11.95%  "
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  -0.12%  "
$ws.Range("D31").Value = "8.16"
$ws.Range("E31").Value = "  +3.50%  "
$ws.Range("E32").Value = "  +0.60%  "
$ws.Range("E33").Value = "  +3.90%  "
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("D35").Value = "3.411.78"
$ws.Range("D36").Value = "23.42"
$ws.Range("E36").Value = "  +2.38%  "
$ws.Range("D37").Value = "5.59"
$ws.Range("E37").Value = "  +5.44%  "
$ws.Range("D38").Value = "6.98"
$ws.Range("E38").Value = "  +3.22%  "
$ws.Range("D39").Value = "1.55"
$ws.Range("E39").Value = "  +3.75%  "
$ws.Range("D40").Value = "162.76"
$ws.Range("E40").Value = "  -0.94%  "
$ws.Range("E41").Value = "  +3.31%  "
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("B43").Value = "ONDO"
$ws.Range("C43").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D43").Value = "1.23"
$ws.Range("E43").Value = "  +13.19%  "
$ws.Range("D44").Value = "4.44"
$ws.Range("E44").Value = "  +3.50%  "
$ws.Range("D45").Value = "41.66"
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("E46").Value = "  -2.48%  "
$ws.Range("D47").Value = "1.62"
$ws.Range("E47").Value = "  +2.85%  "
$ws.Range("D48").Value = "23.50"
$ws.Range("E48").Value = "  +2.41%  "
$ws.Range("D49").Value = "6.92"
$ws.Range("E49").Value = "  +3.62%  "
$ws.Range("D50").Value = "23.22"
$ws.Range("E50").Value = "  +14.57%  "
$ws.Range("D51").Value = "0.901"
$ws.Range("E51").Value = "  +5.32%  "
